# Apply the "Updated symbol list" commit changes to the crypto price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($sheet, [string]$addr, [string]$val)
    $cell = $sheet.Range($addr)
    # Prefix with an apostrophe so numeric-looking strings ("242.49", "0.0005811", ...)
    # are stored as text (preserving exact formatting, e.g. trailing zeros) while
    # keeping the cell's original "General" number format, matching the source data
    # which stores these as inline strings rather than numbers.
    $cell.Value = "'" + $val
}

# --- Column D (Price) numeric-looking text updates ---
Set-TextValue $ws "D2"  "242.49"
Set-TextValue $ws "D4"  "5.222"
Set-TextValue $ws "D5"  "0.05597"
Set-TextValue $ws "D7"  "6.374"
Set-TextValue $ws "D8"  "0.8051"
Set-TextValue $ws "D9"  "0.9343"
Set-TextValue $ws "D10" "0.1426"
Set-TextValue $ws "D11" "0.07281"
Set-TextValue $ws "D12" "0.03124"
Set-TextValue $ws "D14" "0.09282"
Set-TextValue $ws "D15" "3.603"
Set-TextValue $ws "D16" "0.001645"
Set-TextValue $ws "D18" "0.0005811"
Set-TextValue $ws "D19" "0.006347"
Set-TextValue $ws "D20" "0.004980"
Set-TextValue $ws "D23" "0.0003101"
Set-TextValue $ws "D24" "3.762"
Set-TextValue $ws "D25" "2.101"
Set-TextValue $ws "D40" "0.03917"
Set-TextValue $ws "D41" "0.006880"
Set-TextValue $ws "D44" "0.007496"
Set-TextValue $ws "D45" "0.00005943"
Set-TextValue $ws "D47" "0.0005501"
Set-TextValue $ws "D48" "0.6826"
Set-TextValue $ws "D49" "0.07100"

# --- Column E (composite label) text updates ---
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E47").Value = "46ACDXExchangeACXT"

# --- Rows 42 / 43: BKEXToken and CEJI swap places ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.003400"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1031"
$ws.Range("E43").Value = "42BKEXTokenBKK"
